$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 554.8261
$ws.Range("J17").Value = 355.8421
$ws.Range("L17").Value = 1067.5263
$ws.Range("N17").Value = -1403.5263

$ws.Range("H112").Value = 1254.7869
$ws.Range("J112").Value = 1299
$ws.Range("L112").Value = 3897
$ws.Range("N112").Value = -6113

$ws.Range("H131").Value = 4288.933
$ws.Range("I131").Value = 5212.5
$ws.Range("J131").Value = 3953.0908
$ws.Range("K131").Value = 15637.5
$ws.Range("L131").Value = 11859.2724
$ws.Range("M131").Value = -10597.5
$ws.Range("N131").Value = -21939.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 20000
$ws.Range("J9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("N9").Value = -20340

$ws.Range("H20").Value = 20000
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20540

$ws.Range("H37").Value = 25048.545
$ws.Range("J37").Value = 30091.75
$ws.Range("L37").Value = 30091.75
$ws.Range("N37").Value = -30637.75

$ws.Range("H44").Value = 41469.816
$ws.Range("J44").Value = 41469.816
$ws.Range("L44").Value = 41469.816
$ws.Range("N44").Value = -42445.816

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = ""
$ws.Range("N46").Value = 0

$ws.Range("H55").Value = 41584
$ws.Range("J55").Value = 41584
$ws.Range("L55").Value = 41584
$ws.Range("N55").Value = -42214

$ws.Range("H80").Value = 35766.855
$ws.Range("J80").Value = 35766.855
$ws.Range("L80").Value = 35766.855
$ws.Range("N80").Value = -37762.855

$ws.Range("H83").Value = 35766.855
$ws.Range("J83").Value = 35766.855
$ws.Range("L83").Value = 107300.565
$ws.Range("N83").Value = -117284.565

$ws.Range("H102").Value = 3036.6667
$ws.Range("I102").Value = 2555
$ws.Range("K102").Value = 2555
$ws.Range("M102").Value = -933

$ws.Range("H110").Value = 991.75
$ws.Range("I110").Value = 991.75
$ws.Range("K110").Value = 991.75
$ws.Range("M110").Value = 1053.25

$ws.Range("H137").Value = 42000
$ws.Range("J137").Value = 42000
$ws.Range("L137").Value = 42000
$ws.Range("N137").Value = -52200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = ""
$ws.Range("N45").Value = 0

$ws.Range("H59").Value = 43333.332
$ws.Range("J59").Value = 43333.332
$ws.Range("L59").Value = 43333.332
$ws.Range("N59").Value = -45027.332

$ws.Range("H99").Value = 5913.875
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 5913.875
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = ""
$ws.Range("M99").Value = 5913.875
$ws.Range("N99").Value = -8909.875

$ws.Range("H105").Value = 1638.9546
$ws.Range("I105").Value = 1597.9474
$ws.Range("J105").Value = 1898.6666
$ws.Range("K105").Value = 1597.9474
$ws.Range("L105").Value = 1898.6666
$ws.Range("M105").Value = 149.0526
$ws.Range("N105").Value = -5392.6666

$ws.Range("H137").Value = 45760
$ws.Range("J137").Value = 45760
$ws.Range("L137").Value = 45760
$ws.Range("N137").Value = -55960

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 933.3333
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1800
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 1800
$ws.Range("M22").Value = -150
$ws.Range("N22").Value = -2500

$ws.Range("H31").Value = 2778.0908
$ws.Range("I31").Value = 892.44446
$ws.Range("J31").Value = 5040.8667
$ws.Range("K31").Value = 892.44446
$ws.Range("L31").Value = 5040.8667
$ws.Range("M31").Value = -597.44446
$ws.Range("N31").Value = -5630.8667

$ws.Range("H34").Value = 2778.0908
$ws.Range("I34").Value = 892.44446
$ws.Range("J34").Value = 5040.8667
$ws.Range("K34").Value = 892.44446
$ws.Range("L34").Value = 5040.8667
$ws.Range("M34").Value = -690.44446
$ws.Range("N34").Value = -5444.8667

$ws.Range("H99").Value = 6899758.5
$ws.Range("I99").Value = 14287234
$ws.Range("J99").Value = 4781.3335
$ws.Range("K99").Value = 14287234
$ws.Range("L99").Value = 4781.3335
$ws.Range("M99").Value = -14285736
$ws.Range("N99").Value = -7777.3335

$ws.Range("H126").Value = 6899758.5
$ws.Range("I126").Value = 14287234
$ws.Range("J126").Value = 4781.3335
$ws.Range("K126").Value = 42861702
$ws.Range("L126").Value = 14344.0005
$ws.Range("M126").Value = -42859232
$ws.Range("N126").Value = -19284.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5208940.5
$ws.Range("I113").Value = 650.6
$ws.Range("J113").Value = 8929147
$ws.Range("K113").Value = 1951.8
$ws.Range("L113").Value = 26787441
$ws.Range("M113").Value = 218.1999999999998
$ws.Range("N113").Value = -26791781

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 35151.2
$ws.Range("J46").Value = 35151.2
$ws.Range("L46").Value = 35151.2
$ws.Range("N46").Value = -35463.2

$ws.Range("H47").Value = 8000
$ws.Range("J47").Value = 8000
$ws.Range("L47").Value = 8000
$ws.Range("N47").Value = -9136

$ws.Range("H120").Value = 35033.332
$ws.Range("J120").Value = 35033.332
$ws.Range("L120").Value = 35033.332
$ws.Range("N120").Value = -44709.332

$ws.Range("H137").Value = 37135
$ws.Range("J137").Value = 37135
$ws.Range("L137").Value = 37135
$ws.Range("N137").Value = -47335

$ws.Range("H139").Value = 65933.336
$ws.Range("J139").Value = 65933.336
$ws.Range("L139").Value = 65933.336
$ws.Range("N139").Value = -76213.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1980.8462
$ws.Range("I22").Value = 1541.6666
$ws.Range("J22").Value = 2357.2856
$ws.Range("K22").Value = 1541.6666
$ws.Range("L22").Value = 2357.2856
$ws.Range("M22").Value = -1246.6666
$ws.Range("N22").Value = -2947.2856

$ws.Range("H27").Value = 1980.8462
$ws.Range("I27").Value = 1541.6666
$ws.Range("J27").Value = 2357.2856
$ws.Range("K27").Value = 1541.6666
$ws.Range("L27").Value = 2357.2856
$ws.Range("M27").Value = -1434.6666
$ws.Range("N27").Value = -2571.2856

$ws.Range("H40").Value = 8388.4
$ws.Range("I40").Value = 7276.8
$ws.Range("K40").Value = 7276.8
$ws.Range("M40").Value = -7140.8

$ws.Range("H46").Value = 2317.9565
$ws.Range("I46").Value = 1783.4166
$ws.Range("J46").Value = 2901.0908
$ws.Range("K46").Value = 1783.4166
$ws.Range("L46").Value = 2901.0908
$ws.Range("M46").Value = -1595.4166
$ws.Range("N46").Value = -3277.0908

$ws.Range("H50").Value = 49295.8
$ws.Range("J50").Value = 49295.8
$ws.Range("L50").Value = 49295.8
$ws.Range("N50").Value = -50569.8

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = ""
$ws.Range("N133").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 23257
$ws.Range("I37").Value = 7666.3335
$ws.Range("J37").Value = 70029
$ws.Range("K37").Value = 7666.3335
$ws.Range("L37").Value = 70029
$ws.Range("M37").Value = -7463.3335
$ws.Range("N37").Value = -70435

$ws.Range("H46").Value = 65660
$ws.Range("J46").Value = 65660
$ws.Range("L46").Value = 65660
$ws.Range("N46").Value = -66122

$ws.Range("H96").Value = 1872516.6
$ws.Range("I96").Value = 63412.25
$ws.Range("J96").Value = 4767083.5
$ws.Range("K96").Value = 63412.25
$ws.Range("L96").Value = 4767083.5
$ws.Range("M96").Value = -62039.25
$ws.Range("N96").Value = -4769829.5

$ws.Range("H113").Value = 1240
$ws.Range("I113").Value = 980
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 2940
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -770
$ws.Range("N113").Value = -8840

$ws.Range("H122").Value = 2842.353
$ws.Range("I122").Value = 1451.238
$ws.Range("J122").Value = 5089.5386
$ws.Range("K122").Value = 4353.714
$ws.Range("L122").Value = 15268.6158
$ws.Range("M122").Value = -1903.714
$ws.Range("N122").Value = -20168.6158

$ws.Range("H134").Value = 65660
$ws.Range("J134").Value = 65660
$ws.Range("L134").Value = 196980
$ws.Range("N134").Value = -202050

$ws.Range("H138").Value = 37649.668
$ws.Range("J138").Value = 37649.668
$ws.Range("L138").Value = 37649.668
$ws.Range("N138").Value = -47929.668
